$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 154, shifting all existing data (rows 154-215)
# down to rows 157-218. The inserted rows copy the formatting of the row
# above (row 153), which carries the date-number-format style used by
# column D in this table.
$ws.Rows("154:156").Insert()

# Row 154: new "Especial" quote for 2022-11-10 (serial 44875)
$ws.Cells.Item(154, 1).Value  = 8
$ws.Cells.Item(154, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(154, 3).Value  = "Coquimbo"
$ws.Cells.Item(154, 4).Value  = 44875
$ws.Cells.Item(154, 5).Value  = 4
$ws.Cells.Item(154, 6).Value  = "Fruta"
$ws.Cells.Item(154, 7).Value  = 100107
$ws.Cells.Item(154, 8).Value  = "Otros"
$ws.Cells.Item(154, 9).Value  = 100107002
$ws.Cells.Item(154, 10).Value = "Chirimoya"
$ws.Cells.Item(154, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(154, 12).Value = "Especial"
$ws.Cells.Item(154, 13).Value = 240
$ws.Cells.Item(154, 14).Value = 19000
$ws.Cells.Item(154, 15).Value = 20000
$ws.Cells.Item(154, 16).Value = 19500
$ws.Cells.Item(154, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(154, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(154, 19).Value = 1950
$ws.Cells.Item(154, 20).Value = 10

# Row 155: new "Primera" quote for 2022-11-10 (serial 44875)
$ws.Cells.Item(155, 1).Value  = 8
$ws.Cells.Item(155, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(155, 3).Value  = "Coquimbo"
$ws.Cells.Item(155, 4).Value  = 44875
$ws.Cells.Item(155, 5).Value  = 4
$ws.Cells.Item(155, 6).Value  = "Fruta"
$ws.Cells.Item(155, 7).Value  = 100107
$ws.Cells.Item(155, 8).Value  = "Otros"
$ws.Cells.Item(155, 9).Value  = 100107002
$ws.Cells.Item(155, 10).Value = "Chirimoya"
$ws.Cells.Item(155, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(155, 12).Value = "Primera"
$ws.Cells.Item(155, 13).Value = 400
$ws.Cells.Item(155, 14).Value = 16000
$ws.Cells.Item(155, 15).Value = 17000
$ws.Cells.Item(155, 16).Value = 16500
$ws.Cells.Item(155, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(155, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(155, 19).Value = 1650
$ws.Cells.Item(155, 20).Value = 10

# Row 156: new "Segunda" quote for 2022-11-10 (serial 44875)
$ws.Cells.Item(156, 1).Value  = 8
$ws.Cells.Item(156, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(156, 3).Value  = "Coquimbo"
$ws.Cells.Item(156, 4).Value  = 44875
$ws.Cells.Item(156, 5).Value  = 4
$ws.Cells.Item(156, 6).Value  = "Fruta"
$ws.Cells.Item(156, 7).Value  = 100107
$ws.Cells.Item(156, 8).Value  = "Otros"
$ws.Cells.Item(156, 9).Value  = 100107002
$ws.Cells.Item(156, 10).Value = "Chirimoya"
$ws.Cells.Item(156, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(156, 12).Value = "Segunda"
$ws.Cells.Item(156, 13).Value = 240
$ws.Cells.Item(156, 14).Value = 13000
$ws.Cells.Item(156, 15).Value = 14000
$ws.Cells.Item(156, 16).Value = 13500
$ws.Cells.Item(156, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(156, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(156, 19).Value = 1350
$ws.Cells.Item(156, 20).Value = 10
